$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1133
$ws1.Range("F5").Value = 184
$ws1.Range("F14").Value = 12845
$ws1.Range("F16").Value = 5276

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1133
$ws4.Range("F5").Value = 184
$ws4.Range("F14").Value = 12845
$ws4.Range("F18").Value = 5276
